# TD-6649 add Business Concept file manager domain name field
#
# The "domain" header in A1 is renamed to "domain_external_id" (a new
# field). The underlying data row (A2:E2) keeps its original values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "domain_external_id"

# Selection moves to A2 after the edit.
$ws.Range("A2").Select() | Out-Null
